# "Readd venturis to pressure drop"
#
# The sheet used to have two "Venturi" part rows (row 3 for LOx, row 18 for
# Fuel) whose data had been pushed out to orphaned columns L:U (disconnected
# from the main A:J part table) and whose Part Name cells held placeholder
# text ("Part 2" / "Part 22"). This restores them as proper rows in the main
# A:J table, named "LOx Venturi" and "Fuel Venturi".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: LOx Venturi -------------------------------------------------
# Bring the cell formatting for L3:U3 back into A3:J3 (same row, just
# shifted to the left so it lines up under the real headers).
$ws.Range("L3:U3").Copy()
$ws.Range("A3:J3").PasteSpecial(-4122) # xlPasteFormats

# --- Row 18: Fuel Venturi ------------------------------------------------
$ws.Range("L18:U18").Copy()
$ws.Range("A18:J18").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

# Fill in the restored values for the two Venturi rows.
$ws.Range("A3").Value = "LOx Venturi"
$ws.Range("B3").Value = "Venturi"
$ws.Range("H3").Value = 0.75
$ws.Range("I3").Value = 0.065
$ws.Range("J3").Value = 0.001

$ws.Range("A18").Value = "Fuel Venturi"
$ws.Range("B18").Value = "Venturi"
$ws.Range("H18").Value = 0.75
$ws.Range("I18").Value = 0.065
$ws.Range("J18").Value = 0.001

# The old data lived out in columns K:U (disconnected from the table); wipe
# that whole block now that it has been relocated into A:J. This also
# shrinks the sheet's used range back down from A1:U22 to A1:J22.
$ws.Range("K1:U22").Clear()

# Match the saved selection left behind by the edit.
$ws.Range("A19").Select()
